$d = $word.ActiveDocument

# Locate the "Presentación del proyecto" heading; a duplicated, centered
# paragraph containing only an inline picture immediately follows it and
# must be removed entirely (the picture is already shown earlier, under
# the "Introducción" heading).
$findRange = $d.Content
$found = $findRange.Find.Execute("Presentación del proyecto", $false, $false,
                                  $false, $false, $false, $true, 1, $false,
                                  "", 0)

if ($found) {
    $headingEnd = $findRange.End

    # Find the inline picture that appears right after the heading
    # (the nearest InlineShape whose range starts at/after headingEnd).
    $dupShape = $null
    $dupShapeStart = -1
    for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
        $shape = $d.InlineShapes.Item($i)
        $shapeStart = $shape.Range.Start
        if ($shapeStart -ge $headingEnd) {
            if ($dupShape -eq $null -or $shapeStart -lt $dupShapeStart) {
                $dupShape = $shape
                $dupShapeStart = $shapeStart
            }
        }
    }

    if ($dupShape -ne $null) {
        $dupPara = $dupShape.Range.Paragraphs.Item(1)
        $dupRange = $dupPara.Range

        # Safety check: only delete if this paragraph indeed holds the
        # duplicated inline picture (and no other visible text), and it
        # sits immediately after the heading (no other content between).
        if ($dupRange.InlineShapes.Count -gt 0 -and $dupRange.Text.Trim() -eq "" `
            -and $dupRange.Start -lt ($headingEnd + 5)) {
            $dupRange.Delete()
        }
    }
}
